# "Implement separated weather sources." - split the combined
# "Chateau <region>" label in column A into a bare Chateau name, now that
# the region is already carried independently in column B.
#
# Sheet "First Rating" (the active sheet) holds 14 blocks of 25 rows
# (vintages 1994-2018) - one block per Chateau. For every block we
# rewrite column A to the short Chateau name; column B keeps the same
# region text it already had (its shared-string slot just changes
# because the table is being de-duplicated/re-ordered upstream - the
# visible text for B does not change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("First Rating")

$blocks = @(
    @{ Start = 2;   End = 26;  Chateau = "Château Lafite Rothschild";                     Region = "Médoc" },
    @{ Start = 27;  End = 51;  Chateau = "Château Latour";                                 Region = "Médoc" },
    @{ Start = 52;  End = 76;  Chateau = "Château Mouton Rothschild";                      Region = "Médoc" },
    @{ Start = 77;  End = 101; Chateau = "Château Cheval Blanc";                           Region = "Saint Emilion" },
    @{ Start = 102; End = 126; Chateau = "Château Ausone";                                 Region = "Saint Emilion" },
    @{ Start = 127; End = 151; Chateau = "Château Margaux";                                Region = "Médoc" },
    @{ Start = 152; End = 176; Chateau = "Château Haut-Brion";                             Region = "Médoc" },
    @{ Start = 177; End = 201; Chateau = "Château Angélus";                                Region = "Saint Emilion" },
    @{ Start = 202; End = 226; Chateau = "Château Trotte Vieille";                         Region = "Saint Emilion" },
    @{ Start = 227; End = 251; Chateau = "Château Palmer";                                 Region = "Médoc" },
    @{ Start = 252; End = 276; Chateau = "Petrus";                                         Region = "Pomerol" },
    @{ Start = 277; End = 301; Chateau = "Château Pavie";                                  Region = "Saint Emilion" },
    @{ Start = 302; End = 326; Chateau = "Château Léoville Las Cases";                     Region = "Médoc" },
    @{ Start = 327; End = 351; Chateau = "Château Pichon-Longueville Comtesse de Lalande";  Region = "Médoc" }
)

foreach ($block in $blocks) {
    for ($row = $block.Start; $row -le $block.End; $row++) {
        $ws.Cells.Item($row, 1).Value = $block.Chateau
        $ws.Cells.Item($row, 2).Value = $block.Region
    }
}

# Restore the workbook/window chrome that Excel re-wrote on save: the
# sheet was scrolled down and the selection moved to A31.
$ws.Range("A31").Select()
